$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D19: remove "treatment satisfaction [OGMS:0000090]; " and the space before the BOM
$ws.Range("D19").Value = "injury [OGMS:0000102];" + [char]0xFEFF + "disease [OGMS:0000031]"

# Add new row 28 for GSSO import
$ws.Range("A28").Value = "GSSO"
$ws.Range("B28").Value = "http://purl.obolibrary.org/obo/gsso.owl"
$ws.Range("C28").Value = "entity [BFO:0000001]"
$ws.Range("D28").Value = "advocacy organisation [GSSO:005379]; health organisation [GSSO:007328]; human rights organisation [GSSO:003501]; non-profit organisation [GSSO:004615]; money [GSSO:010609]"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "all"
$ws.Range("F28").Value = ""
$ws.Range("F28").Font.Bold = $false
